{"js": "// Fix \"Saal\" -> \"Saalbau\" in the \"\u00dcbrige Lokale...\" paragraph, and remove the\n// internal review question paragraph (\"Frage Hanni: ...\") together with the\n// blank paragraph that immediately followed it.\n\nconst body = context.document.body;\n\n// 1) \"...IBA Geb\u00e4ude, Saal zwei Mal, KIFF...\" -> \"...IBA Geb\u00e4ude, Saalbau zwei Mal, KIFF...\"\nconst hits = body.search(\"Saal zwei Mal, KIFF\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < hits.items.length; i++) {\n  hits.items[i].insertText(\"Saalbau zwei Mal, KIFF\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Remove the \"Frage Hanni\" question paragraph and the empty paragraph after it.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Frage Hanni\") !== -1) {\n    // Remove the trailing blank paragraph first so indices stay valid.\n    if (i + 1 < items.length) {\n      items[i + 1].delete();\n    }\n    items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"...IBA Gebaude, Saal zwei Mal, KIFF...\" -> \"...IBA Gebaude, Saalbau zwei Mal, KIFF...\"\n$range = $d.Content\n$found = $range.Find.Execute(\n    \"Saal zwei Mal, KIFF\",  # FindText\n    $false,                 # MatchCase\n    $false,                 # MatchWholeWord\n    $false,                 # MatchWildcards\n    $false,                 # MatchSoundsLike\n    $false,                 # MatchAllWordForms\n    $true,                  # Forward\n    1,                      # Wrap (wdFindContinue)\n    $false,                 # Format\n    \"Saalbau zwei Mal, KIFF\", # ReplaceWith\n    2                       # Replace (wdReplaceAll)\n)\n\n# 2) Remove the \"Frage Hanni\" question paragraph and the blank paragraph\n#    immediately following it (walk backwards so indices stay valid).\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Frage Hanni*\") {\n        if ($i + 1 -le $d.Paragraphs.Count) {\n            $d.Paragraphs.Item($i + 1).Range.Delete()\n        }\n        $p.Range.Delete()\n        break\n    }\n}\n"}
